# Applies the "Add files via upload" edit to the Path to Graduation workbook.
# Source: replaces the placeholder student name/ID, shuffles the Fall
# 2022 / Fall 2023 / Fall 2024 course blocks to their updated course codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: student name / id -------------------------------------------
# E1 is a numeric-looking ID that must stay text (leading zero). Force the
# text number format before assigning it, then drop back to the default
# "Normal" style so we don't leave a stray custom format behind.
$ws.Range("C1").Value = "Billy Bob"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "0321472904323"
$ws.Range("E1").Style = "Normal"

# --- Fall 2022 / Summer 2022 block (rows 6-8) -----------------------------
# Old: C6=CYBR 3106/D6=3 ; C7=CPSC 1302/D7=3
# New: C6=CPSC 1302/D6=3 ; A7=CYBR 3106/B7=3 ; C7=CPSC 2108/D7=3 ; C8=CYBR 3108/D8=3
$ws.Cells.Item(6, 3).Value = "CPSC 1302"
$ws.Cells.Item(6, 4).Value = 3

$ws.Cells.Item(7, 1).Value = "CYBR 3106"
$ws.Cells.Item(7, 2).Value = 3
$ws.Cells.Item(7, 3).Value = "CPSC 2108"
$ws.Cells.Item(7, 4).Value = 3

$ws.Cells.Item(8, 3).Value = "CYBR 3108"
$ws.Cells.Item(8, 4).Value = 3

# --- Fall 2023 / Spring 2023 block (rows 13-17) ---------------------------
$ws.Cells.Item(13, 1).Value = "CPSC 4155"
$ws.Cells.Item(13, 2).Value = 3
$ws.Cells.Item(13, 3).Value = "CPSC 4135"
$ws.Cells.Item(13, 4).Value = 3

$ws.Cells.Item(14, 1).Value = "DSCI 3111"
$ws.Cells.Item(14, 2).Value = 3
$ws.Cells.Item(14, 3).Value = "CPSC 4175"
$ws.Cells.Item(14, 4).Value = 3

$ws.Cells.Item(15, 1).Value = "CPSC 3165"
$ws.Cells.Item(15, 2).Value = 3
$ws.Cells.Item(15, 3).Value = "CPSC 6180"
$ws.Cells.Item(15, 4).Value = 3

$ws.Cells.Item(16, 1).Value = "CPSC 4111"
$ws.Cells.Item(16, 2).Value = 3
$ws.Cells.Item(16, 3).Value = "CPSC 6185"
$ws.Cells.Item(16, 4).Value = 3

$ws.Cells.Item(17, 1).Value = "CPSC 4148"
$ws.Cells.Item(17, 2).Value = 3

# --- Fall 2024 / Spring 2024 block (rows 22-24) ---------------------------
# Old row22: A=CPSC 4175/B=3, C=CPSC 6985/D=4
# New row22: A=CPSC 6985/B=4  (C/D cleared)
$ws.Cells.Item(22, 1).Value = "CPSC 6985"
$ws.Cells.Item(22, 2).Value = 4
$ws.Range("C22:D22").ClearContents()

# Old row23: A=CPSC 6180/B=3, C=CPSC 4000/D=0
# New row23: A=CPSC 4000/B=0  (C/D cleared)
$ws.Cells.Item(23, 1).Value = "CPSC 4000"
$ws.Cells.Item(23, 2).Value = 0
$ws.Range("C23:D23").ClearContents()

# Old row24: A=CPSC 6185/B=3 -- removed entirely
$ws.Range("A24:F24").ClearContents()
